$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this shifts the existing rows 16-32 down to 17-33
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record
$ws.Range("A16").Value2 = 10
$ws.Range("B16").Value2 = "Vega Modelo de Temuco"
$ws.Range("C16").Value2 = "La Araucanía"
$ws.Range("D16").Value2 = 44771
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E16").Value2 = 9
$ws.Range("F16").Value2 = 100112042
$ws.Range("G16").Value2 = "Locoto"
$ws.Range("H16").Value2 = "Sin especificar"
$ws.Range("I16").Value2 = "Primera"
$ws.Range("J16").Value2 = 30
$ws.Range("K16").Value2 = 3300
$ws.Range("L16").Value2 = 3300
$ws.Range("M16").Value2 = 3300
$ws.Range("N16").Value2 = "$/kilo"
$ws.Range("O16").Value2 = "Región de Arica y Parinacota"
$ws.Range("P16").Value2 = 3300
$ws.Range("Q16").Value2 = 1
$ws.Range("R16").Value2 = "Hortaliza"
